$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2644.7
$ws.Range("J40").Value = 2640.9412
$ws.Range("L40").Value = 2640.9412
$ws.Range("N40").Value = -2990.9412
$ws.Range("H43").Value = 1654.2727
$ws.Range("I43").Value = 1299.75
$ws.Range("J43").Value = 1856.8572
$ws.Range("K43").Value = 1299.75
$ws.Range("L43").Value = 1856.8572
$ws.Range("M43").Value = -1230.75
$ws.Range("N43").Value = -1994.8572
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H87").Value = 67485
$ws.Range("J87").Value = 92098.75
$ws.Range("L87").Value = 92098.75
$ws.Range("N87").Value = -94594.75
$ws.Range("H90").Value = 67485
$ws.Range("J90").Value = 92098.75
$ws.Range("L90").Value = 276296.25
$ws.Range("N90").Value = -288776.25
$ws.Range("H101").Value = 2487.125
$ws.Range("J101").Value = 4985
$ws.Range("L101").Value = 14955
$ws.Range("N101").Value = -18199
$ws.Range("H112").Value = 1584.4386
$ws.Range("J112").Value = 1672.7451
$ws.Range("L112").Value = 5018.2353
$ws.Range("N112").Value = -7234.2353
$ws.Range("H138").Value = 3368.8684
$ws.Range("I138").Value = 2195.3684
$ws.Range("J138").Value = 4542.3687
$ws.Range("K138").Value = 6586.1052
$ws.Range("L138").Value = 13627.1061
$ws.Range("M138").Value = -1446.1052
$ws.Range("N138").Value = -23907.1061
$ws.Range("H141").Value = 1361.625
$ws.Range("I141").Value = 1385
$ws.Range("J141").Value = 1198
$ws.Range("K141").Value = 4155
$ws.Range("L141").Value = 3594
$ws.Range("M141").Value = 1025
$ws.Range("N141").Value = -13954
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1963883.8
$ws.Range("I61").Value = 2384208
$ws.Range("J61").Value = 2370
$ws.Range("K61").Value = 2384208
$ws.Range("L61").Value = 2370
$ws.Range("M61").Value = -2383996
$ws.Range("N61").Value = -2794
$ws.Range("H74").Value = 4810883.5
$ws.Range("I74").Value = 6252307
$ws.Range("J74").Value = 6140
$ws.Range("K74").Value = 6252307
$ws.Range("L74").Value = 6140
$ws.Range("M74").Value = -6251433
$ws.Range("N74").Value = -7888
$ws.Range("H77").Value = 4810883.5
$ws.Range("I77").Value = 6252307
$ws.Range("J77").Value = 6140
$ws.Range("K77").Value = 31261535
$ws.Range("L77").Value = 30700
$ws.Range("M77").Value = -31257167
$ws.Range("N77").Value = -39436
$ws.Range("H122").Value = 4999.7144
$ws.Range("I122").Value = 4999
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 14997
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -12547
$ws.Range("N122").Value = -19900
$ws.Range("H132").Value = 467758.88
$ws.Range("I132").Value = 568811.0600000001
$ws.Range("K132").Value = 1706433.18
$ws.Range("M132").Value = -1703903.18
$ws.Range("H134").Value = 75499.664
$ws.Range("J134").Value = 75499.664
$ws.Range("L134").Value = 75499.664
$ws.Range("N134").Value = -85639.664
$ws.Range("H136").Value = 1963883.8
$ws.Range("I136").Value = 2384208
$ws.Range("J136").Value = 2370
$ws.Range("K136").Value = 7152624
$ws.Range("L136").Value = 7110
$ws.Range("M136").Value = -7150074
$ws.Range("N136").Value = -12210
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 99951
$ws.Range("J135").Value = 99951
$ws.Range("L135").Value = 99951
$ws.Range("N135").Value = -110091
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 45000
$ws.Range("I17").Value = 57000
$ws.Range("K17").Value = 57000
$ws.Range("M17").Value = -56826
$ws.Range("H62").Value = 999
$ws.Range("J62").Value = 999
$ws.Range("L62").Value = 999
$ws.Range("N62").Value = -2247
$ws.Range("H65").Value = 999
$ws.Range("J65").Value = 999
$ws.Range("L65").Value = 4995
$ws.Range("N65").Value = -11235
$ws.Range("H68").Value = 72352.81
$ws.Range("J68").Value = 72352.81
$ws.Range("L68").Value = 72352.81
$ws.Range("N68").Value = -73850.81
$ws.Range("H71").Value = 72352.81
$ws.Range("J71").Value = 72352.81
$ws.Range("L71").Value = 217058.43
$ws.Range("N71").Value = -224546.43
$ws.Range("H99").Value = 5457.1113
$ws.Range("I99").Value = 4375.7144
$ws.Range("K99").Value = 4375.7144
$ws.Range("M99").Value = -2877.7144
$ws.Range("H105").Value = 44081.75
$ws.Range("I105").Value = 69926.39999999999
$ws.Range("J105").Value = 1007.3333
$ws.Range("K105").Value = 69926.39999999999
$ws.Range("L105").Value = 1007.3333
$ws.Range("M105").Value = -68179.39999999999
$ws.Range("N105").Value = -4501.3333
$ws.Range("H122").Value = 3629.75
$ws.Range("I122").Value = 1914
$ws.Range("J122").Value = 5345.5
$ws.Range("K122").Value = 5742
$ws.Range("L122").Value = 16036.5
$ws.Range("M122").Value = -3292
$ws.Range("N122").Value = -20936.5
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H126").Value = 5457.1113
$ws.Range("I126").Value = 4375.7144
$ws.Range("K126").Value = 13127.1432
$ws.Range("M126").Value = -10657.1432
$ws.Range("H140").Value = 90483.2
$ws.Range("J140").Value = 90483.2
$ws.Range("L140").Value = 90483.2
$ws.Range("N140").Value = -100843.2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 182441.55
$ws.Range("I14").Value = 182441.55
$ws.Range("K14").Value = 547324.6499999999
$ws.Range("M14").Value = -547151.6499999999
$ws.Range("H37").Value = 85772.55
$ws.Range("J37").Value = 85772.55
$ws.Range("L37").Value = 257317.65
$ws.Range("N37").Value = -257541.65
$ws.Range("H120").Value = 26265.75
$ws.Range("I120").Value = 15015
$ws.Range("J120").Value = 37516.5
$ws.Range("K120").Value = 45045
$ws.Range("L120").Value = 112549.5
$ws.Range("M120").Value = -40207
$ws.Range("N120").Value = -122225.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1675.2894
$ws.Range("I102").Value = 1207.2667
$ws.Range("K102").Value = 1207.2667
$ws.Range("M102").Value = 414.7333000000001
$ws.Range("H132").Value = 710687.5
$ws.Range("I132").Value = 754730.5
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 2264191.5
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -2261661.5
$ws.Range("N132").Value = -23060
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H130").Value = 79729.89
$ws.Range("J130").Value = 79729.89
$ws.Range("L130").Value = 79729.89
$ws.Range("N130").Value = -89769.89
$ws.Range("H136").Value = 8865.5
$ws.Range("I136").Value = 8630.799999999999
$ws.Range("J136").Value = 9452.25
$ws.Range("K136").Value = 25892.4
$ws.Range("L136").Value = 28356.75
$ws.Range("M136").Value = -23342.4
$ws.Range("N136").Value = -33456.75
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 15000
$ws.Range("I7").Value = 5000
$ws.Range("K7").Value = 5000
$ws.Range("M7").Value = -4887
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("H132").Value = 6291829.5
$ws.Range("I132").Value = 6710751.5
$ws.Range("J132").Value = 8000
$ws.Range("K132").Value = 20132254.5
$ws.Range("L132").Value = 24000
$ws.Range("M132").Value = -20129724.5
$ws.Range("N132").Value = -29060
